$wb = $excel.ActiveWorkbook

# Add new test-step info to the "Test 1" sheet first (keeps shared-string order
# lined up with the authored workbook)
$test1 = $wb.Worksheets.Item("Test 1")
$test1.Range("B2").Value = "open eclipse"
$test1.Range("C2").Value = "eclipse opens"
$test1.Range("C14").Select()

# Insert a brand new worksheet before "Test 1" and rename it "Test Ideas"
$ideas = $wb.Worksheets.Add($wb.Worksheets.Item("Test 1"))
$ideas.Name = "Test Ideas"

# Column A is wide, single column of test-idea notes
$ideas.Columns.Item(1).ColumnWidth = 76

$ideas.Range("A2").Value = "Test for mouse click of button to progress through slides"
$ideas.Range("A1").Value = "Test for audio played (different file types)"
$ideas.Range("A3").Value = "Test for text objects dispayed on slide (different fonts + sizes + added features)"
$ideas.Range("A4").Value = "Test for video objects played (also ALL control functions/buttons)"

$ideas.Activate()
$ideas.Range("A4").Select()
